# Apply the NCR_HEALTH.xlsx header restructuring described in the commit:
#  - Uppercase a couple of header labels
#  - Remove the placeholder "-" values from columns I and L (data rows)
#  - Split the old "AA" column into five new "No. of Sites ..." header
#    columns (AA-AE) that reuse the bold/bordered header style, and push the
#    "Status as of ..." column out to AF
#  - Move the dropdown data validation from AA2:AA17 to AF2:AF17

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text tweaks -----------------------------------------------
$ws.Range("I1").Value = "TOTAL PHYSICAL TARGET"
$ws.Range("L1").Value = "BATCH"

# --- Remove placeholder "-" cells in columns I and L (rows 2-17) ------
$ws.Range("I2:I17").ClearContents()
$ws.Range("L2:L17").ClearContents()

# --- Insert the five new header cells (AA1:AE1), reusing the existing --
# --- bold/boxed header formatting (same style as the rest of row 1) ---
$ws.Range("Z1").Copy()
$ws.Range("AA1").PasteSpecial(-4122)
$ws.Range("AA1").Value = "No. of Sites Reverted"

$ws.Range("Z1").Copy()
$ws.Range("AB1").PasteSpecial(-4122)
$ws.Range("AB1").Value = "No. of Sites Not yet started"

$ws.Range("Z1").Copy()
$ws.Range("AC1").PasteSpecial(-4122)
$ws.Range("AC1").Value = "No. of Sites Under Procurement"

$ws.Range("Z1").Copy()
$ws.Range("AD1").PasteSpecial(-4122)
$ws.Range("AD1").Value = "No. of Sites On Going"

$ws.Range("Z1").Copy()
$ws.Range("AE1").PasteSpecial(-4122)
$ws.Range("AE1").Value = "No. of Sites Completed"

# --- Move the "Status as of ..." header out to AF1 (unstyled, like the --
# --- original AA1) -----------------------------------------------------
$ws.Range("AF1").Value = "Status as of July 4, 2025"

# --- Move the dropdown validation from AA2:AA17 to AF2:AF17 ------------
$ws.Range("AA2:AA17").Validation.Delete()
$v = $ws.Range("AF2:AF17").Validation
$v.Add(3, 1, 1, "=DropdownOptions!`$A`$1:`$A`$7")
$v.IgnoreBlank = $true
$v.InCellDropdown = $true
$v.ShowInput = $false
$v.ShowError = $false
